# "changing FALSE to False" -- column I (rows 2-41) held a boolean formula
# =FALSE() that rendered as "FALSE" via a custom "TRUE"/"FALSE" number
# format. The author replaced those with the literal text string "False".
#
# Repro approach:
#   1. Put a text-producing formula ("False") in the range, which avoids
#      Excel's automatic bool-literal coercion that a plain
#      Range.Value = "False" assignment would trigger.
#   2. Copy / Paste-Special values-only over itself to freeze the formula
#      result down to a literal (shared-string) value, t="s".
#   3. Re-apply a Text number format ("@") so the cells stop using the old
#      "TRUE"/"FALSE" custom format (numFmtId 164 -> 49).
#   4. A few extra blank rows (42-45) picked up the same style while the
#      user was doing this (visible in the dimension/selection growing to
#      row 45), so stamp the same number format on I42:I45 too.
#   5. Leave the selection where the author's last selection ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$boolRange = $ws.Range("I2:I41")

$boolRange.Formula = '="False"'
$boolRange.Copy()
$boolRange.PasteSpecial(-4163)  # xlPasteValues
$ws.Application.CutCopyMode = 0

$ws.Range("I2:I45").NumberFormat = "@"

$ws.Range("H42:J46").Select()
